$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2263.0454
$ws.Range("J112").Value = 2610.389
$ws.Range("L112").Value = 7831.167
$ws.Range("N112").Value = -10047.167
$ws.Range("H116").Value = 2416.3635
$ws.Range("I116").Value = 1611.8572
$ws.Range("J116").Value = 3824.25
$ws.Range("K116").Value = 1611.8572
$ws.Range("L116").Value = 3824.25
$ws.Range("M116").Value = 1830.1428
$ws.Range("N116").Value = -10708.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1498.81
$ws.Range("I32").Value = 1405.3298
$ws.Range("J32").Value = 2963.3333
$ws.Range("K32").Value = 1405.3298
$ws.Range("L32").Value = 2963.3333
$ws.Range("M32").Value = -1118.3298
$ws.Range("N32").Value = -3537.3333
$ws.Range("H61").Value = 142858380
$ws.Range("I61").Value = 166667620
$ws.Range("K61").Value = 166667620
$ws.Range("M61").Value = -166667408
$ws.Range("H74").Value = 1269.3334
$ws.Range("I74").Value = 962.2353000000001
$ws.Range("J74").Value = 2574.5
$ws.Range("K74").Value = 962.2353000000001
$ws.Range("L74").Value = 2574.5
$ws.Range("M74").Value = -88.23530000000005
$ws.Range("N74").Value = -4322.5
$ws.Range("H77").Value = 1269.3334
$ws.Range("I77").Value = 962.2353000000001
$ws.Range("J77").Value = 2574.5
$ws.Range("K77").Value = 4811.1765
$ws.Range("L77").Value = 12872.5
$ws.Range("M77").Value = -443.1765000000005
$ws.Range("N77").Value = -21608.5
$ws.Range("H132").Value = 1828.772
$ws.Range("I132").Value = 1408.6842
$ws.Range("K132").Value = 4226.0526
$ws.Range("M132").Value = -1696.0526
$ws.Range("H136").Value = 142858380
$ws.Range("I136").Value = 166667620
$ws.Range("K136").Value = 500002860
$ws.Range("M136").Value = -500000310

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1112.44
$ws.Range("I134").Value = 895.9048
$ws.Range("K134").Value = 2687.7144
$ws.Range("M134").Value = -152.7143999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1622
$ws.Range("I31").Value = 1443.9429
$ws.Range("J31").Value = 2660.6667
$ws.Range("K31").Value = 1443.9429
$ws.Range("L31").Value = 2660.6667
$ws.Range("M31").Value = -1148.9429
$ws.Range("N31").Value = -3250.6667
$ws.Range("H34").Value = 1622
$ws.Range("I34").Value = 1443.9429
$ws.Range("J34").Value = 2660.6667
$ws.Range("K34").Value = 1443.9429
$ws.Range("L34").Value = 2660.6667
$ws.Range("M34").Value = -1241.9429
$ws.Range("N34").Value = -3064.6667
$ws.Range("H58").Value = 1508.8823
$ws.Range("I58").Value = 1058.8823
$ws.Range("J58").Value = 1958.8823
$ws.Range("K58").Value = 1058.8823
$ws.Range("L58").Value = 1958.8823
$ws.Range("M58").Value = -855.8823
$ws.Range("N58").Value = -2364.8823
$ws.Range("H105").Value = 982.75
$ws.Range("I105").Value = 955
$ws.Range("K105").Value = 955
$ws.Range("M105").Value = 792
$ws.Range("H132").Value = 1291.9688
$ws.Range("I132").Value = 863.7826
$ws.Range("K132").Value = 2591.3478
$ws.Range("M132").Value = -61.34780000000001
$ws.Range("H136").Value = 1508.8823
$ws.Range("I136").Value = 1058.8823
$ws.Range("J136").Value = 1958.8823
$ws.Range("K136").Value = 3176.6469
$ws.Range("L136").Value = 5876.6469
$ws.Range("M136").Value = -626.6468999999997
$ws.Range("N136").Value = -10976.6469
$ws.Range("H137").Value = 70996
$ws.Range("J137").Value = 77994.664
$ws.Range("L137").Value = 77994.664
$ws.Range("N137").Value = -88194.664
$ws.Range("H141").Value = 683953.3
$ws.Range("J141").Value = 683953.3
$ws.Range("L141").Value = 683953.3
$ws.Range("N141").Value = -694313.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 985.55554
$ws.Range("I97").Value = 792.5
$ws.Range("J97").Value = 1140
$ws.Range("K97").Value = 2377.5
$ws.Range("L97").Value = 3420
$ws.Range("M97").Value = -1881.5
$ws.Range("N97").Value = -4412
$ws.Range("H102").Value = 3169.4285
$ws.Range("J102").Value = 3169.4285
$ws.Range("L102").Value = 9508.2855
$ws.Range("N102").Value = -14376.2855
$ws.Range("H140").Value = 2611.537
$ws.Range("I140").Value = 1863
$ws.Range("J140").Value = 3479.84
$ws.Range("K140").Value = 5589
$ws.Range("L140").Value = 10439.52
$ws.Range("M140").Value = -409
$ws.Range("N140").Value = -20799.52

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 432.6087
$ws.Range("I107").Value = 329.6
$ws.Range("J107").Value = 625.75
$ws.Range("K107").Value = 329.6
$ws.Range("L107").Value = 625.75
$ws.Range("M107").Value = 1590.4
$ws.Range("N107").Value = -4465.75
$ws.Range("H126").Value = 1642.1177
$ws.Range("I126").Value = 1438.2858
$ws.Range("J126").Value = 2593.3333
$ws.Range("K126").Value = 4314.857400000001
$ws.Range("L126").Value = 7779.999899999999
$ws.Range("M126").Value = -1844.857400000001
$ws.Range("N126").Value = -12719.9999
$ws.Range("H132").Value = 2181.818
$ws.Range("I132").Value = 1935.6129
$ws.Range("J132").Value = 2768.923
$ws.Range("K132").Value = 5806.8387
$ws.Range("L132").Value = 8306.769
$ws.Range("M132").Value = -3276.8387
$ws.Range("N132").Value = -13366.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 2980
$ws.Range("I29").Value = 2980
$ws.Range("K29").Value = 2980
$ws.Range("M29").Value = -2685
$ws.Range("H31").Value = 4155
$ws.Range("I31").Value = 2000
$ws.Range("J31").Value = 4586
$ws.Range("K31").Value = 2000
$ws.Range("L31").Value = 4586
$ws.Range("M31").Value = -1752
$ws.Range("N31").Value = -5082
$ws.Range("H32").Value = 2708.3333
$ws.Range("I32").Value = 1650
$ws.Range("K32").Value = 1650
$ws.Range("M32").Value = -1333
$ws.Range("H34").Value = 2500
$ws.Range("J34").Value = 3000
$ws.Range("L34").Value = 3000
$ws.Range("N34").Value = -3344
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H136").Value = 1833.3334
$ws.Range("I136").Value = 1500
$ws.Range("K136").Value = 4500
$ws.Range("M136").Value = -1950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 13159339
$ws.Range("I122").Value = 14707267
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 44121801
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -44119351
$ws.Range("N122").Value = -10750
$ws.Range("H126").Value = 76928400
$ws.Range("I126").Value = 111114460
$ws.Range("K126").Value = 333343380
$ws.Range("M126").Value = -333340910
$ws.Range("H133").Value = 34903.75
$ws.Range("J133").Value = 34903.75
$ws.Range("L133").Value = 34903.75
$ws.Range("N133").Value = -45023.75
$ws.Range("H136").Value = 1665.7059
$ws.Range("I136").Value = 1345.3334
$ws.Range("J136").Value = 2026.125
$ws.Range("K136").Value = 4036.0002
$ws.Range("L136").Value = 6078.375
$ws.Range("M136").Value = -1486.0002
$ws.Range("N136").Value = -11178.375
